# "added some initial vf implementation"
#
# The Th/U/Pu (thorium/uranium/plutonium) related rows and their
# gas/liquid-suffixed species labels are removed from the condensation
# sequence tables, and Table 1's remaining species labels are normalized
# from their "_g"/"_liq" suffixed forms to the plain element/oxide names.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Table 1: drop the Th/U/Pu reaction rows (10-12) and rename the
# remaining Product 1 / Product 2 / Reactant species labels so they no
# longer carry the "_g" / "_liq" phase suffixes.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Table 1")

$ws1.Range("A10:H12").EntireRow.Delete()

$product1 = @("Si", "Mg", "Fe", "Ca", "Al", "Ti", "Na", "K")
for ($i = 0; $i -lt $product1.Length; $i++) {
    $ws1.Cells.Item(2 + $i, 2).Value = $product1[$i]
}

$ws1.Range("C2:C9").Value = "O"

$reactant = @("SiO2", "MgO", "FeO", "CaO", "Al2O3", "TiO2", "Na2O", "K2O")
for ($i = 0; $i -lt $reactant.Length; $i++) {
    $ws1.Cells.Item(2 + $i, 6).Value = $reactant[$i]
}

# ---------------------------------------------------------------------
# Table 2: drop the Th/U/Pu element rows (28-38).
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Table 2")
$ws2.Range("A28:P38").EntireRow.Delete()

# ---------------------------------------------------------------------
# View/selection bookkeeping to mirror the saved UI state: Table 1 is no
# longer the active tab, Table 2's frozen pane keeps scrolled to where
# the deleted rows used to be, and Table 3 becomes the active tab with a
# portrait page setup.
# ---------------------------------------------------------------------
$ws1.Range("A10:H12").Select()

$ws2.Activate()
$ws2.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $false
$excel.ActiveWindow.FreezePanes = $true
$ws2.Range("A28:R38").Select()

$ws3 = $wb.Worksheets.Item("Table 3")
$ws3.PageSetup.Orientation = 1
$ws3.Activate()
$ws3.Range("C37").Select()
